# StyleTagTemplate.xlsx update:
#  - add a new worksheet "class only" at the end of the workbook, containing
#    five jt:style / class-attribute example rows (becomes the active/selected
#    tab, matching Excel's behavior of tracking the last-inserted / active sheet)
#  - the previously-selected first sheet ("alignment") loses its tabSelected flag
#    automatically once a different sheet becomes active/selected

$wb = $excel.ActiveWorkbook

# --- add the new sheet after the last existing sheet -----------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "class only"

# --- populate the new sheet -------------------------------------------------
# Values are entered in the same order they first appear in the template so
# that the resulting shared-string table is built up in that same order.
$newSheet.Range("B2").Value  = '<jt:style class="redBoxCenter">redBoxCenter</jt:style>'
$newSheet.Range("B6").Value  = '<jt:style class="doesNotExist">doesNotExist</jt:style>'
$newSheet.Range("B4").Value  = '<jt:style class="blueBoldBigText">blueBoldBigText</jt:style>'
$newSheet.Range("B8").Value  = '<jt:style class="blueBoldBigText" style="font-color: green; font-weight: normal; font-italic: true">blueBoldBigText with green non-bold italic</jt:style>'
$newSheet.Range("B10").Value = '<jt:style class="redBoxCenter; blueBoldBigText">redBoxCenter; blueBoldBigText</jt:style>'

# Column B is widened, same as the other template sheets.
$newSheet.Columns("B").ColumnWidth = 35.67

# Make the new sheet the active / selected tab.
$newSheet.Activate()
$newSheet.Select()

Write-Output ("Worksheets: " + $wb.Worksheets.Count)
Write-Output ("Active sheet: " + $wb.ActiveSheet.Name)
